# "added common utility for runmodes"
#
# TestCases sheet: OpenAccountTest's Runmode flips from N -> Y (run every
# test by default now).
#
# TestData sheet: the AddCustomerTest data block gains two more duplicate
# data rows (now runs for manish + jyoti twice -- N becomes Y for jyoti's
# original row), and the OpenAccountTest data block also gains two more
# duplicate rows (re-using the same two customers), with the jyoti row
# there flipping from Y to N.

$wb = $excel.ActiveWorkbook
$wsCases = $wb.Worksheets.Item("TestCases")
$wsData  = $wb.Worksheets.Item("TestData")

# --- TestCases sheet: OpenAccountTest runmode N -> Y -----------------
$wsCases.Range("B3").Value2 = "Y"

# --- TestData sheet -----------------------------------------------------
# Current layout (before edit):
#   1  AddCustomerTest
#   2  Runmode | firstname | lastname | postcode
#   3  Y | manish | k | P6767
#   4  N | jyoti  | k | X7878
#   5  (blank)
#   6  OpenAccountTest
#   7  Runmode | customer  | currency
#   8  Y | manish k | Rupee
#   9  Y | jyoti k  | Dollar

# 1) Flip row 4 (jyoti / AddCustomerTest) runmode from N to Y.
$wsData.Range("A4").Value2 = "Y"

# 2) Insert two new rows right after row 4, pushing the blank separator
#    row and the whole OpenAccountTest block down by two rows.
$wsData.Rows("5:6").Insert()

# 3) Populate the two freshly inserted rows as duplicates of the
#    (now all-"Y") AddCustomerTest data rows.
$wsData.Range("A5:D5").Value2 = $wsData.Range("A3:D3").Value2
$wsData.Range("A6:D6").Value2 = $wsData.Range("A4:D4").Value2

# After the insert, the OpenAccountTest block now lives at:
#   8  OpenAccountTest
#   9  Runmode | customer  | currency
#   10 Y | manish k | Rupee
#   11 Y | jyoti k  | Dollar   (still Y at this point)

# 4) Flip the jyoti OpenAccountTest runmode from Y to N.
$wsData.Range("A11").Value2 = "N"

# 5) Append two more duplicate rows (manish/jyoti) at the bottom of the
#    OpenAccountTest block.
$wsData.Range("A12:C12").Value2 = $wsData.Range("A10:C10").Value2
$wsData.Range("A13:C13").Value2 = $wsData.Range("A11:C11").Value2

# --- View / selection state ---------------------------------------------
# Final state: TestData's selection rests on the last appended row, and
# TestCases becomes the active (tab-selected) sheet with B2 selected.
$wsData.Select()
$wsData.Range("A13").Select()

$wsCases.Select()
$wsCases.Range("B2").Select()

